$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2024-10-31 12:39:51"
$ws.Range("B6").Value = "Success"

$ws.Range("A7").Value = "2024-10-31 12:42:59"
$ws.Range("B7").Value = "Success"

$ws.Range("A8").Value = "2024-10-31 12:47:45"
$ws.Range("B8").Value = "Success"
